$d = $word.ActiveDocument

# --- 1) Insert a new Heading1 paragraph before the very first paragraph ---
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphBefore() | Out-Null
$newHeading = $d.Paragraphs.Item(1)
$newHeading.Range.Text = "Redo (Albedo for only EP different(can you see a difference) )"
$newHeading.Style = "Heading 1"

# --- 2) Merge the split [UWG]/[IDF] runs in the "Bypass" row of the first table ---
# (Table 1, row 4 "Bypass", column 2 still has the text split across two runs per
#  paragraph; column 3 in the same row is already a single run per paragraph.)
$tbl = $d.Tables.Item(1)
$cell = $tbl.Cell(4, 2)
$cellStart = $cell.Range.Start

# Paragraph 1: "[UWG]Albedo = 0.15"
$p1End = $cellStart + 18
$p1 = $d.Range($cellStart, $p1End)
$p1.Text = "[UWG]Albedo = 0.15 "
$p1b = $d.Range($cellStart, $cellStart + 19)
$p1b.Text = "[UWG]Albedo = 0.15"

# Paragraph 2: "[IDF]Solar Absorptance = 0.5"
$p2Start = $cellStart + 19
$p2End = $cellStart + 47
$p2 = $d.Range($p2Start, $p2End)
$p2.Text = "[IDF]Solar Absorptance = 0.5 "
$p2b = $d.Range($p2Start, $p2Start + 29)
$p2b.Text = "[IDF]Solar Absorptance = 0.5"

# --- 3) Highlight six CVRMSE "Bypass"-column figures in yellow ---
$targets = @("7.71", "9.45", "7.49", "7.68", "11.19", "7.56")
foreach ($t in $targets) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($t, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.HighlightColorIndex = 7
    }
}
